$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- Row 25: flip sign of F, G, H ---
$ws.Range("F25").Value = 33009.148504422461
$ws.Range("G25").Value = 44391.613505947455
$ws.Range("H25").Value = 134313.08701799484

# --- Rows 66-70 and 72-78: flip sign of F, G, H (row 71 untouched, values are 0) ---
$ws.Range("F66").Value = -49513.722756633695
$ws.Range("G66").Value = -66587.420258921178
$ws.Range("H66").Value = -201469.63052699226

$ws.Range("F67").Value = -29708.23365398022
$ws.Range("G67").Value = -39952.452155352716
$ws.Range("H67").Value = -120881.77831619537

$ws.Range("F68").Value = -9902.7445513267448
$ws.Range("G68").Value = -13317.484051784249
$ws.Range("H68").Value = -40293.926105398481

$ws.Range("F69").Value = 39610.978205306958
$ws.Range("G69").Value = 53269.936207136947
$ws.Range("H69").Value = 161175.70442159381

$ws.Range("F70").Value = 19805.489102653482
$ws.Range("G70").Value = 26634.968103568481
$ws.Range("H70").Value = 80587.852210796918

$ws.Range("F72").Value = -16504.57425221123
$ws.Range("G72").Value = -22195.806752973727
$ws.Range("H72").Value = -67156.54350899742

$ws.Range("F73").Value = 19247.475701909207
$ws.Range("G73").Value = 25884.536288774452
$ws.Range("H73").Value = 78317.314925009865

$ws.Range("F74").Value = -19247.475701909207
$ws.Range("G74").Value = -25884.536288774452
$ws.Range("H74").Value = -78317.314925009865

$ws.Range("F75").Value = 19247.475701909207
$ws.Range("G75").Value = 25884.536288774452
$ws.Range("H75").Value = 78317.314925009865

$ws.Range("F76").Value = -19247.475701909207
$ws.Range("G76").Value = -25884.536288774452
$ws.Range("H76").Value = -78317.314925009865

$ws.Range("F77").Value = 19247.475701909207
$ws.Range("G77").Value = 25884.536288774452
$ws.Range("H77").Value = 78317.314925009865

$ws.Range("F78").Value = -16504.57425221123
$ws.Range("G78").Value = -22195.806752973727
$ws.Range("H78").Value = -67156.54350899742

# --- Rows 101-104: fill in previously empty F, G, H ---
$ws.Range("F101").Value = 17.385849761672155
$ws.Range("G101").Value = 23.380970369145317
$ws.Range("H101").Value = 70.742423168183237

$ws.Range("F102").Value = -30.014999999999997
$ws.Range("G102").Value = -40.365000000000002
$ws.Range("H102").Value = -122.12999999999997

$ws.Range("F103").Value = 8.6929248808360793
$ws.Range("G103").Value = 11.69048518457266
$ws.Range("H103").Value = 35.371211584091625

$ws.Range("F104").Value = 30.014999999999997
$ws.Range("G104").Value = 40.365000000000002
$ws.Range("H104").Value = 122.12999999999997

# --- Rows 110-113: fill in previously empty F ---
$ws.Range("F110").Value = 1.5
$ws.Range("F111").Value = 1.2
$ws.Range("F112").Value = 1.5
$ws.Range("F113").Value = 1.2

# --- Rows 125-128: fill in previously empty F, G, H ---
$ws.Range("F125").Value = 26.078774642508233
$ws.Range("G125").Value = 35.071455553717968
$ws.Range("H125").Value = 106.11363475227486

$ws.Range("F126").Value = -36.017999999999994
$ws.Range("G126").Value = -48.438000000000002
$ws.Range("H126").Value = -146.55599999999998

$ws.Range("F127").Value = 13.03938732125412
$ws.Range("G127").Value = 17.535727776858987
$ws.Range("H127").Value = 53.056817376137438

$ws.Range("F128").Value = 36.017999999999994
$ws.Range("G128").Value = 48.438000000000002
$ws.Range("H128").Value = 146.55599999999998

# --- Rows 135-137: fill in previously empty F ---
$ws.Range("F135").Value = 3258.721877189379
$ws.Range("F136").Value = 33.582666144355763
$ws.Range("F137").Value = 32.750118334874607

# --- Sheet view: scroll/selection position ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 106
$win.ScrollColumn = 1
$ws.Range("H136").Select()
